# Update "想去人数" (F) and "最低票价" (G) figures on both the "展览"
# and "全部类型" sheets, matching the refreshed scrape output.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("F2").Value = 203
    $ws.Range("F3").Value = 435
    $ws.Range("F4").Value = 12781
    # Row 5 sold out: numeric min-price of 1 becomes the text "已售罄"
    $ws.Range("G5").Value = "已售罄"
    $ws.Range("F6").Value = 174
    $ws.Range("F7").Value = 35
    $ws.Range("F10").Value = 214
    $ws.Range("F11").Value = 461
    $ws.Range("F15").Value = 45
    $ws.Range("F16").Value = 390
    $ws.Range("F17").Value = 5466
    $ws.Range("F19").Value = 32
    $ws.Range("F20").Value = 951
    $ws.Range("F23").Value = 102
}
